$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$newRows = @(
    @{ idx = 263; date = "2016-08-09"; open = 33541.26; low = 33350.38; high = 33596.61; close = 33528.52 },
    @{ idx = 264; date = "2016-08-10"; open = 33666.24; low = 33249.26; high = 33822.85; close = 33504.81 },
    @{ idx = 265; date = "2016-08-11"; open = 33466.81; low = 33218.79; high = 33586.21; close = 33278.06 }
)

$r = 265
foreach ($row in $newRows) {
    $ws.Cells.Item($r, 1).Value = $row.idx
    $ws.Cells.Item($r, 2).Value = "HSP"
    $ws.Cells.Item($r, 3).NumberFormat = "@"
    $ws.Cells.Item($r, 3).Value = $row.date
    $ws.Cells.Item($r, 4).Value = $row.open
    $ws.Cells.Item($r, 5).Value = $row.low
    $ws.Cells.Item($r, 6).Value = $row.high
    $ws.Cells.Item($r, 7).Value = $row.close
    $r = $r + 1
}
